$d = $word.ActiveDocument

$d.Content.Find.Execute("939÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "640÷7=", 2) | Out-Null
$d.Content.Find.Execute("689÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "297÷8=", 2) | Out-Null
$d.Content.Find.Execute("561÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "399÷9=", 2) | Out-Null
$d.Content.Find.Execute("821÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "290÷5=", 2) | Out-Null
$d.Content.Find.Execute("556÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "732÷9=", 2) | Out-Null
$d.Content.Find.Execute("184÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "542÷6=", 2) | Out-Null
$d.Content.Find.Execute("965÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "254÷3=", 2) | Out-Null
$d.Content.Find.Execute("998÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "912÷5=", 2) | Out-Null
$d.Content.Find.Execute("269÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "103÷8=", 2) | Out-Null
$d.Content.Find.Execute("175÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "132÷3=", 2) | Out-Null
$d.Content.Find.Execute("377÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "654÷9=", 2) | Out-Null
$d.Content.Find.Execute("973÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "755÷2=", 2) | Out-Null
$d.Content.Find.Execute("170÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "638÷2=", 2) | Out-Null
$d.Content.Find.Execute("514÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "594÷7=", 2) | Out-Null
$d.Content.Find.Execute("476÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "885÷7=", 2) | Out-Null
$d.Content.Find.Execute("571÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "171÷2=", 2) | Out-Null
$d.Content.Find.Execute("130÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "335÷3=", 2) | Out-Null
$d.Content.Find.Execute("839÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "814÷9=", 2) | Out-Null
$d.Content.Find.Execute("718÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "769÷2=", 2) | Out-Null
$d.Content.Find.Execute("403÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "512÷7=", 2) | Out-Null
$d.Content.Find.Execute("824÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "537÷3=", 2) | Out-Null
$d.Content.Find.Execute("583÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "216÷9=", 2) | Out-Null
$d.Content.Find.Execute("714÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "634÷7=", 2) | Out-Null
$d.Content.Find.Execute("913÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "700÷8=", 2) | Out-Null
$d.Content.Find.Execute("215÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "737÷7=", 2) | Out-Null
